$d = $word.ActiveDocument

# --- Change 1: "The goal of this project..." paragraph ---
$p3 = $d.Paragraphs(3)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End)
$rng3.Text = "The goal of this project is to provide updated information on unemployment, stock prices, and the number of cases and deaths in the United States by state and week. "

# --- Change 2: merge "Our solution..." paragraph with the following
#     "First of all..." paragraph, then replace with new combined text ---
$p5 = $d.Paragraphs(5)
$markPos = $p5.Range.End - 1
$d.Range($markPos, $markPos + 1).Delete()

$p5merged = $d.Paragraphs(5)
$rng5 = $d.Range($p5merged.Range.Start, $p5merged.Range.End)
$rng5.Text = "Our solution will be useful to different kinds of users. Firstly, this information could be used by investors, should there be another pandemic, to help predict when the stock crash would happen. Secondly, it could be used by economic analysts outside of buying stocks to help predict costs due to unemployment or otherwise understand the potential impact of a pandemic. Thirdly, city, state, and government officials could use this information to help understand growth rate of the virus, understand how many supplies they may need, or see where in the country the virus is spreading. This could help with supply chain management. "

# --- Change 3: remove the "second type of user", "for example", the blank
#     paragraph between them, and the "Finally, ..." paragraphs entirely ---
$p7 = $d.Paragraphs(7)
$p10 = $d.Paragraphs(10)
$delRng = $d.Range($p7.Range.Start, $p10.Range.End)
$delRng.Delete()
